$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column E to (as close as the engine allows to) fit the new, longer value
$ws.Columns.Item(5).ColumnWidth = 10

# Copy the formatting (styles) of row 3's date / IsShortSell cells onto row 6's
# matching cells before writing values, so the new row reuses the existing
# cellXfs entries (s="1") instead of minting new ones.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null

$ws.Range("G3").Copy() | Out-Null
$ws.Range("G6").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# New row 6 data (a new trade record)
$ws.Range("A6").Value = 42649.655729166669

$ws.Range("B6").Value = $false

$ws.Range("C6").Value = 9888.5400000000009
$ws.Range("D6").Value = 9931.74
$ws.Range("E6").Value = 313.26998900000001
$ws.Range("F6").Value = 310.52999999999997

$ws.Range("G6").Value = $false

$ws.Range("H6").Value = -0.87

$ws.Range("I6").Value = $false
